# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion note text with new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.05 = 27853.37 pesos`n✅ 27853.37 pesos = 7.02 = 972.94 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update N10/O10/N12/O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 141.85
$wsTasas.Range("O10").Value = 3951
$wsTasas.Range("N12").Value = 3965
$wsTasas.Range("O12").Value = 138.5
